$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D15").Value = 26
$wsForecast.Range("D17").Value = 25

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "502"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "25"

$wsSummary.Range("B15").NumberFormat = "@"
$wsSummary.Range("B15").Value = "2025-04-20"
